# Update inclass data after 1_intro finished
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the date (2019-01-22 -> serial 43487) to A2:A6, styled like a date
$ws.Range("A2:A6").Value = 43487
$ws.Range("A2:A6").NumberFormat = "d-mmm"

# Add note to F2
$ws.Range("F2").Value = "See 1_intro.Rmd"

# Update minute values
$ws.Range("E5").Value = 25
$ws.Range("E6").Value = 5

# Update the selection / view
$ws.Range("F5").Select()
